{"js": "// Update the multiplication problems in the table to the new values\n// described by the diff. Each \"find\" string is a unique literal that\n// occurs exactly once in the document body, so a simple search +\n// replace per pair is safe and order independent.\nconst replacements = [\n  { find: \"51\u00d727=\", replace: \"73\u00d787=\" },\n  { find: \"14\u00d711=\", replace: \"37\u00d729=\" },\n  { find: \"99\u00d783=\", replace: \"34\u00d794=\" },\n  { find: \"69\u00d711=\", replace: \"47\u00d738=\" },\n  { find: \"75\u00d752=\", replace: \"30\u00d792=\" },\n  { find: \"99\u00d765=\", replace: \"17\u00d766=\" },\n  { find: \"51\u00d724=\", replace: \"85\u00d730=\" },\n  { find: \"89\u00d756=\", replace: \"80\u00d772=\" },\n  { find: \"88\u00d740=\", replace: \"69\u00d733=\" },\n  { find: \"34\u00d739=\", replace: \"76\u00d717=\" },\n  { find: \"34\u00d798=\", replace: \"67\u00d723=\" },\n  { find: \"19\u00d752=\", replace: \"14\u00d745=\" },\n  { find: \"55\u00d793=\", replace: \"73\u00d727=\" },\n  { find: \"77\u00d783=\", replace: \"51\u00d795=\" },\n  { find: \"34\u00d711=\", replace: \"51\u00d761=\" },\n  { find: \"32\u00d788=\", replace: \"79\u00d753=\" },\n  { find: \"98\u00d797=\", replace: \"25\u00d798=\" },\n  { find: \"64\u00d759=\", replace: \"43\u00d789=\" },\n  { find: \"52\u00d734=\", replace: \"68\u00d792=\" },\n  { find: \"96\u00d746=\", replace: \"81\u00d791=\" },\n  { find: \"96\u00d793=\", replace: \"31\u00d714=\" },\n  { find: \"48\u00d730=\", replace: \"46\u00d725=\" },\n  { find: \"39\u00d753=\", replace: \"40\u00d760=\" },\n  { find: \"29\u00d775=\", replace: \"61\u00d716=\" },\n  { find: \"63\u00d781=\", replace: \"66\u00d717=\" },\n];\n\nconst body = context.document.body;\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(replace, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the multiplication problems in the table to the new values\n# described by the diff. Each \"find\" string is a unique literal that\n# occurs exactly once in the document, so Find/Replace per pair (with\n# a fresh Find range from $d.Content each time) is safe and order\n# independent.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"51\u00d727=\", \"73\u00d787=\"),\n  @(\"14\u00d711=\", \"37\u00d729=\"),\n  @(\"99\u00d783=\", \"34\u00d794=\"),\n  @(\"69\u00d711=\", \"47\u00d738=\"),\n  @(\"75\u00d752=\", \"30\u00d792=\"),\n  @(\"99\u00d765=\", \"17\u00d766=\"),\n  @(\"51\u00d724=\", \"85\u00d730=\"),\n  @(\"89\u00d756=\", \"80\u00d772=\"),\n  @(\"88\u00d740=\", \"69\u00d733=\"),\n  @(\"34\u00d739=\", \"76\u00d717=\"),\n  @(\"34\u00d798=\", \"67\u00d723=\"),\n  @(\"19\u00d752=\", \"14\u00d745=\"),\n  @(\"55\u00d793=\", \"73\u00d727=\"),\n  @(\"77\u00d783=\", \"51\u00d795=\"),\n  @(\"34\u00d711=\", \"51\u00d761=\"),\n  @(\"32\u00d788=\", \"79\u00d753=\"),\n  @(\"98\u00d797=\", \"25\u00d798=\"),\n  @(\"64\u00d759=\", \"43\u00d789=\"),\n  @(\"52\u00d734=\", \"68\u00d792=\"),\n  @(\"96\u00d746=\", \"81\u00d791=\"),\n  @(\"96\u00d793=\", \"31\u00d714=\"),\n  @(\"48\u00d730=\", \"46\u00d725=\"),\n  @(\"39\u00d753=\", \"40\u00d760=\"),\n  @(\"29\u00d775=\", \"61\u00d716=\"),\n  @(\"63\u00d781=\", \"66\u00d717=\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n  $find = $d.Content.Find\n  $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
